$d = $word.ActiveDocument

# 1) Replace the opening sentence.
$d.Content.Find.Execute(
    "Este trabajo es para aprender a utilizar el ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Estoy aprendiendo a utilizar ", 2)

# 2) Replace the closing clause (after the "git" run).
$d.Content.Find.Execute(
    " y el GITHUB y familiarizarnos con los comandos .", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " y GITHUB con los comandos.", 2)

# 3) Append a new, empty trailing paragraph at the end of the document.
$d.Paragraphs.Add() | Out-Null
